$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.242.72"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "2.424.38"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'564.52"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "'144.73"
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").Value = "2.422.40"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "'0.110"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'26.11"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "'0.0000178"
$ws.Range("E15").Value = "  +5.02%  "
$ws.Range("D16").Value = "2.861.95"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("D17").Value = "61.962.30"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "2.424.09"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "'11.32"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("D20").Value = "'325.43"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "'4.20"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'65.55"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").Value = "'1.73"
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("D26").Value = "'9.05"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'588.83"
$ws.Range("E27").Value = "  +13.14%  "
$ws.Range("D28").Value = "2.539.48"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0950"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.46"
$ws.Range("E31").Value = "  +4.41%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.26"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "'0.151"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "'5.75"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'4.81"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").Value = "'154.10"
$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("D40").Value = "'0.384"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").Value = "'18.73"
$ws.Range("D42").Value = "'1.84"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'2.36"
$ws.Range("E44").Value = "  +7.79%  "
$ws.Range("D45").Value = "'150.25"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").Value = "'3.66"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").Value = "'0.0540"
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("D48").Value = "'20.49"
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").Value = "'0.593"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").Value = "'0.0926"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("E51").Value = "  +1.64%  "
